$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4345088601112366
$ws.Range("B1").Value = 0.6414564847946167
$ws.Range("C1").Value = 1.336283683776855
$ws.Range("D1").Value = 5.880020618438721
$ws.Range("E1").Value = 2.33031177520752
